$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the text in E8 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Make E8 the selected / active cell (as reflected by the new <selection> element)
$ws.Range("E8").Select()
